# Apply portfolio data refresh (cached Yahoo scrape) + add Ticker column
# to the Portfolio Summary sheet.

$wb = $excel.ActiveWorkbook

$wsLog = $wb.Worksheets.Item("Stock log")
$wsSummary = $wb.Worksheets.Item("Portfolio Summary")

# ---------------------------------------------------------------------
# "Stock log" sheet: refreshed market price / market value / gains for
# the three holdings (two AY lots + SCHD).
# ---------------------------------------------------------------------

# Row 2 - AY, 687 shares
$wsLog.Range("E2").Value = 21.7
$wsLog.Range("I2").Value = 14907.9
$wsLog.Range("J2").Value = -7.86
$wsLog.Range("L2").Value = -4.24

# Row 3 - AY, 400 shares
$wsLog.Range("E3").Value = 21.7
$wsLog.Range("I3").Value = 8680
$wsLog.Range("J3").Value = -12
$wsLog.Range("L3").Value = -8.71

# Row 4 - SCHD, 983 shares
$wsLog.Range("E4").Value = 76.23
$wsLog.Range("I4").Value = 74934.09
$wsLog.Range("J4").Value = 0.77
$wsLog.Range("L4").Value = 2.67

# ---------------------------------------------------------------------
# "Portfolio Summary" sheet: insert a new leading "Ticker" column
# (shifting every existing column one place to the right) and populate
# it, then refresh the recalculated totals that moved along with it.
# ---------------------------------------------------------------------

$wsSummary.Columns.Item(1).Insert()

# Give the new column A the same bold/bordered header style used by the
# rest of row 1 (and by the ticker data cells in rows 2-3).
$wsSummary.Range("B1").Copy()
$wsSummary.Range("A1:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsSummary.Range("A1").Value = "Ticker"
$wsSummary.Range("A2").Value = "AY"
$wsSummary.Range("A3").Value = "SCHD"

# Updated Market Value / Capital Gains % / Total Return % (now in columns
# D, G, H after the insert) to match the refreshed "Stock log" figures.
$wsSummary.Range("D2").Value = 23587.9
$wsSummary.Range("G2").Value = -9.43
$wsSummary.Range("H2").Value = -5.94

$wsSummary.Range("D3").Value = 74934.09
$wsSummary.Range("G3").Value = 0.77
$wsSummary.Range("H3").Value = 2.67
